$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Crop Size column (G) to the new crop size (224) for rows 5, 6, 11, 12
$ws.Range("G5").Value = 224
$ws.Range("G6").Value = 224
$ws.Range("G11").Value = 224
$ws.Range("G12").Value = 224

# Update the active selection on the sheet to H15
$ws.Range("H15").Select()
